# table_1_78.xlsx - fix D35:D38 formulas.
# They were "=D3/3" and are being changed to "=($A$1*1.1)/2" (same pattern
# used by D3 itself, just halved), which also updates the dependent E, K
# and the C40 total via recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-write the formula for D35:D38 as one block so Excel stores it as a
# single shared formula (matching how the author authored/typed it once
# and filled down/selected D35:D38).
$ws.Range("D35:D38").Formula = "=(`$A`$1*1.1)/2"

# Reflect the author's new viewport/selection on the sheet.
$ws.Activate()
$win = $excel.ActiveWindow
$ws.Range("D35:D38").Select()
$win.ScrollRow = 16
$win.ScrollColumn = 1
